$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.076.64'
$ws.Range("E2").Value = '  -2.57%  '
$ws.Range("D3").Value = '1.637.48'
$ws.Range("E3").Value = '  -2.57%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''308.57'
$ws.Range("E5").Value = '  -2.02%  '
$ws.Range("D6").Value = '''1.004'
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("D7").Value = '''0.3933'
$ws.Range("E7").Value = '  +0.55%  '
$ws.Range("D8").Value = '''0.3855'
$ws.Range("E8").Value = '  -2.45%  '
$ws.Range("D9").Value = '''1.005'
$ws.Range("E9").Value = '  +0.20%  '
$ws.Range("D10").Value = '''49.97'
$ws.Range("E10").Value = '  -3.48%  '
$ws.Range("D11").Value = '''1.363'
$ws.Range("E11").Value = '  -2.57%  '
$ws.Range("D12").Value = '''0.08551'
$ws.Range("E12").Value = '  -1.22%  '
$ws.Range("D13").Value = '''23.56'
$ws.Range("E13").Value = '  -6.72%  '
$ws.Range("D14").Value = '''7.065'
$ws.Range("E14").Value = '  -3.56%  '
$ws.Range("E15").Value = '  -2.43%  '
$ws.Range("D16").Value = '''7.484'
$ws.Range("E16").Value = '  -3.50%  '
$ws.Range("D17").Value = '1.649.38'
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").Value = '''93.83'
$ws.Range("E18").Value = '  +0.32%  '
$ws.Range("D19").Value = '''0.06918'
$ws.Range("E19").Value = '  -2.33%  '
$ws.Range("D20").Value = '''20.29'
$ws.Range("E20").Value = '  +0.15%  '
$ws.Range("D21").Value = '''6.901'
$ws.Range("E21").Value = '  -2.26%  '
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").Value = '''13.58'
$ws.Range("E23").Value = '  -2.44%  '
$ws.Range("D24").Value = '24.076.87'
$ws.Range("E24").Value = '  -2.53%  '
$ws.Range("D25").Value = '''2.404'
$ws.Range("E25").Value = '  +2.47%  '
$ws.Range("D26").Value = '''2.857'
$ws.Range("E26").Value = '  +0.83%  '
$ws.Range("E27").Value = '  -5.05%  '
$ws.Range("D28").Value = '''157.64'
$ws.Range("E28").Value = '  -3.03%  '
$ws.Range("D29").Value = '''139.86'
$ws.Range("E29").Value = '  -4.21%  '
$ws.Range("D30").Value = '''8.096'
$ws.Range("E30").Value = '  +2.88%  '
$ws.Range("D31").Value = '''5.258'
$ws.Range("E31").Value = '  -10.01%  '
$ws.Range("D32").Value = '''2.480'
$ws.Range("E32").Value = '  +4.27%  '
$ws.Range("D33").Value = '1.824.87'
$ws.Range("E33").Value = '  -0.37%  '
$ws.Range("D34").Value = '''0.08052'
$ws.Range("E34").Value = '  -4.19%  '
$ws.Range("D35").Value = '''6.701'
$ws.Range("E35").Value = '  -3.72%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '''0.02895'
$ws.Range("E36").Value = '  -5.11%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '''0.9659'
$ws.Range("E37").Value = '  -3.34%  '
$ws.Range("D38").Value = '''0.2684'
$ws.Range("E38").Value = '  -4.00%  '
$ws.Range("D39").Value = '''0.09252'
$ws.Range("E39").Value = '  -1.96%  '
$ws.Range("D40").Value = '''10.24'
$ws.Range("E40").Value = '  -3.33%  '
$ws.Range("D41").Value = '''1.423'
$ws.Range("E41").Value = '  -7.68%  '
$ws.Range("D42").Value = '''0.7491'
$ws.Range("E42").Value = '  -5.37%  '
$ws.Range("E43").Value = '  -3.57%  '
$ws.Range("D44").Value = '''16.23'
$ws.Range("E44").Value = '  -2.30%  '
$ws.Range("D45").Value = '''0.6891'
$ws.Range("E45").Value = '  -3.52%  '
$ws.Range("D46").Value = '''2.452'
$ws.Range("E46").Value = '  -4.39%  '
$ws.Range("D47").Value = '''4.095'
$ws.Range("E47").Value = '  -2.06%  '
$ws.Range("E48").Value = '  +0.05%  '
$ws.Range("D49").Value = '''0.08346'
$ws.Range("E49").Value = '  -3.72%  '
$ws.Range("D50").Value = '''1.259'
$ws.Range("E50").Value = '  -6.09%  '
$ws.Range("D51").Value = '''133.13'
$ws.Range("E51").Value = '  -3.43%  '
